$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 47
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.Value = "'2025/10/02"
$dateCell.ClearFormats()
$ws.Cells.Item($row, 2).Value = "木"
$ws.Cells.Item($row, 3).Value = 0
$ws.Cells.Item($row, 4).Value = 22
